# Card16: add a new service "Event" log entry and its related columns.
# This mirrors the workflow of the underlying report generator, which
# re-exports the whole sheet (filling previously blank tracking cells
# with the literal placeholder "nan") whenever a brand new column group
# ("Event" / "Correction" / "Servised by") or a new log row is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# 1) Add the three new header cells (row 1), using the same bold/bordered/
#    centered header formatting as the rest of row 1 (style "s=1").
$ws.Cells.Item(1, 13).Value = "Event"
$ws.Cells.Item(1, 14).Value = "Correction"
$ws.Cells.Item(1, 15).Value = "Servised by"

$header = $ws.Range("M1:O1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# 2) Backfill every currently-blank tracking cell (columns D..L, rows 2..12)
#    with the literal text "nan" -- matching cells that already contain a
#    real value (checkmarks, dates, tonnage figures) are left untouched.
for ($r = 2; $r -le 12; $r++) {
    for ($c = 4; $c -le 12; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value() -eq "") {
            $cell.Value = "nan"
        }
    }
}

# 3) Append the new event row (row 13) describing the latest service record.
#    Column A holds the card number as text (matching A2:A12), so force a
#    text number format before assigning it -- otherwise Excel would infer
#    a numeric value.
$row13 = 13
$ws.Cells.Item($row13, 1).NumberFormat = "@"
$ws.Cells.Item($row13, 1).Value = "16"

$ws.Cells.Item($row13, 12).Value = "1\1\2024"
$ws.Cells.Item($row13, 13).Value = "سلك هالك"
$ws.Cells.Item($row13, 14).Value = "تم تغير سلك"
$ws.Cells.Item($row13, 15).Value = "م.رشدي"
